# Add new team rows to BD_Times (sheet1) and corresponding match rows to
# BD_Jogo (sheet2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BD_Times")
$ws2 = $wb.Worksheets.Item("BD_Jogo")

# --- Sheet1 (BD_Times): append rows 506-513, columns A:I ---
$times = @(
    @("Avai",            1, 0, 1, 0, 0, 2, 7, 4),
    @("Atletico-GO",     0, 1, 0, 0, 2, 0, 4, 7),
    @("ABC",              1, 1, 1, 1, 1, 1, 4, 3),
    @("Sampaio Correia", 0, 1, 1, 1, 1, 1, 3, 4),
    @("Ceara",            1, 1, 0, 0, 1, 0, 8, 4),
    @("Criciuma",         0, 0, 1, 0, 0, 1, 4, 8),
    @("Guarani",          1, 1, 0, 0, 1, 0, 6, 1),
    @("Ponte Preta",      0, 0, 1, 0, 0, 1, 1, 6)
)

$startRow1 = 506
for ($i = 0; $i -lt $times.Count; $i++) {
    $r = $startRow1 + $i
    $row = $times[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[8]
}

# --- Sheet2 (BD_Jogo): append rows 254-257, columns A:E ---
$jogos = @(
    @(0, 2, 11, "Avai",        "Atletico-GO"),
    @(1, 2, 7,  "ABC",          "Sampaio Correia"),
    @(0, 1, 12, "Ceara",        "Criciuma"),
    @(0, 1, 7,  "Guarani",      "Ponte Preta")
)

$startRow2 = 254
for ($i = 0; $i -lt $jogos.Count; $i++) {
    $r = $startRow2 + $i
    $row = $jogos[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
}
